$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: 3 -> 4
$ws.Range("B2").Value = 4

# Delete row 4 entirely (A4=2, B4=1) so dimension becomes A1:B3
$ws.Rows(4).Delete()
